$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A77").Value = '''''AGE_DIFFERENCE'','
$ws.Range("A77").Style = "Normal"
$ws.Range("A78").Value = '''''CHILDREN_IN_HH'','
$ws.Range("A78").Style = "Normal"
$ws.Range("A79").Value = '''''DISTANCEMOVED_10MI'','
$ws.Range("A79").Style = "Normal"
$ws.Range("A80").Value = '''''GENDER_ATTRACTION'','
$ws.Range("A80").Style = "Normal"
$ws.Range("A81").Value = '''''HHINC'','
$ws.Range("A81").Style = "Normal"
$ws.Range("A82").Value = '''''HOW_LONG_AGO_FIRST_COHAB'','
$ws.Range("A82").Style = "Normal"
$ws.Range("A83").Value = '''''HOW_LONG_AGO_FIRST_MET'','
$ws.Range("A83").Style = "Normal"
$ws.Range("A84").Value = '''''HOW_LONG_AGO_FIRST_ROMANTIC'','
$ws.Range("A84").Style = "Normal"
$ws.Range("A85").Value = '''''HOW_LONG_RELATIONSHIP'','
$ws.Range("A85").Style = "Normal"
$ws.Range("A86").Value = '''''PARTNER_MOM_YRSED'','
$ws.Range("A86").Style = "Normal"
$ws.Range("A87").Value = '''''PARTNER_YRSED'','
$ws.Range("A87").Style = "Normal"
$ws.Range("A88").Value = '''''PPAGECAT'','
$ws.Range("A88").Style = "Normal"
$ws.Range("A89").Value = '''''PPHOUSE'','
$ws.Range("A89").Style = "Normal"
$ws.Range("A90").Value = '''''PPHOUSEHOLDSIZE'','
$ws.Range("A90").Style = "Normal"
$ws.Range("A91").Value = '''''Q21A'','
$ws.Range("A91").Style = "Normal"
$ws.Range("A92").Value = '''''Q21B'','
$ws.Range("A92").Style = "Normal"
$ws.Range("A93").Value = '''''Q21C'','
$ws.Range("A93").Style = "Normal"
$ws.Range("A94").Value = '''''Q9'','
$ws.Range("A94").Style = "Normal"
$ws.Range("A95").Value = '''''RELATIONSHIP_QUALITY'','
$ws.Range("A95").Style = "Normal"
$ws.Range("A96").Value = '''''RESPONDENT_MOM_YRSED'','
$ws.Range("A96").Style = "Normal"
$ws.Range("A97").Value = '''''RESPONDENT_YRSED'','
$ws.Range("A97").Style = "Normal"
$ws.Range("A98").Value = '''''ZPFORBORN_CAT'','
$ws.Range("A98").Style = "Normal"
$ws.Range("A99").Value = '''''ZPNHBLACK_CAT'','
$ws.Range("A99").Style = "Normal"
$ws.Range("A100").Value = '''''ZPNHWHITE_CAT'']'
$ws.Range("A100").Style = "Normal"
$ws.Range("A101").Value = 'CORESIDENT'','
$ws.Range("A102").Value = '''''GLBSTATUS'','
$ws.Range("A102").Style = "Normal"
$ws.Range("A103").Value = '''''MARRIED'','
$ws.Range("A103").Style = "Normal"
$ws.Range("A104").Value = '''''MET_THROUGH_AS_COWORKERS'','
$ws.Range("A104").Style = "Normal"
$ws.Range("A105").Value = '''''MET_THROUGH_AS_NEIGHBORS'','
$ws.Range("A105").Style = "Normal"
$ws.Range("A106").Value = '''''MET_THROUGH_FAMILY'','
$ws.Range("A106").Style = "Normal"
$ws.Range("A107").Value = '''''MET_THROUGH_FRIENDS'','
$ws.Range("A107").Style = "Normal"
$ws.Range("A108").Value = '''''PPHHHEAD'','
$ws.Range("A108").Style = "Normal"
$ws.Range("A109").Value = '''''PPMSACAT'','
$ws.Range("A109").Style = "Normal"
$ws.Range("A110").Value = '''''PPNET'','
$ws.Range("A110").Style = "Normal"
$ws.Range("A111").Value = '''''Q31_1'', #refused is -1'
$ws.Range("A111").Style = "Normal"
$ws.Range("A112").Value = '''''Q31_2'',#refused is -1'
$ws.Range("A112").Style = "Normal"
$ws.Range("A113").Value = '''''Q31_3'',#refused is -1'
$ws.Range("A113").Style = "Normal"
$ws.Range("A114").Value = '''''Q31_4'',#refused is -1'
$ws.Range("A114").Style = "Normal"
$ws.Range("A115").Value = '''''Q31_5'',#refused is -1'
$ws.Range("A115").Style = "Normal"
$ws.Range("A116").Value = '''''Q31_6'',#refused is -1'
$ws.Range("A116").Style = "Normal"
$ws.Range("A117").Value = '''''Q31_7'',#refused is -1'
$ws.Range("A117").Style = "Normal"
$ws.Range("A118").Value = '''''Q31_8'',#refused is -1'
$ws.Range("A118").Style = "Normal"
$ws.Range("A119").Value = '''''SAME_SEX_COUPLE'','
$ws.Range("A119").Style = "Normal"
$ws.Range("A120").Value = '''''US_RAISED'','
$ws.Range("A120").Style = "Normal"
$ws.Range("A121").Value = '''''ZPRURAL_CAT'','
$ws.Range("A121").Style = "Normal"
$ws.Range("A122").Value = '''''PARENTAL_APPROVAL'','
$ws.Range("A122").Style = "Normal"
$ws.Range("A123").Value = '''''Q33_1'',  #refused is -1'
$ws.Range("A123").Style = "Normal"
$ws.Range("A124").Value = '''''Q33_2'',#refused is -1'
$ws.Range("A124").Style = "Normal"
$ws.Range("A125").Value = '''''Q33_3'', #refused is -1'
$ws.Range("A125").Style = "Normal"
$ws.Range("A126").Value = '''''Q33_4'', #refused is -1'
$ws.Range("A126").Style = "Normal"
$ws.Range("A127").Value = '''''Q33_5'', #refused is -1'
$ws.Range("A127").Style = "Normal"
$ws.Range("A128").Value = '''''Q33_6'', #refused is -1'
$ws.Range("A128").Style = "Normal"
$ws.Range("A129").Value = '''''Q33_7'',#refused is -1'
$ws.Range("A129").Style = "Normal"
$ws.Range("A130").Value = '''''EITHER_INTERNET_ADJUSTED'
$ws.Range("A130").Style = "Normal"
$ws.Range("A131").Value = 'PAPEVANGELICAL'', #change to 0/1 from 1/2, no refused'
$ws.Range("A132").Value = '''''Q13A'', #change to 0/1 from 1/2, refused is -1'
$ws.Range("A132").Style = "Normal"
$ws.Range("A136").Value = '''''Q28'', #change to 0/1 from 1/2, refused is -1'
$ws.Range("A136").Style = "Normal"
$ws.Range("A138").Value = '''''Q8A'''
$ws.Range("A138").Style = "Normal"
$ws.Range("A137").Value = '''Q7A'''
$ws.Range("A135").Value = '''Q27'', '
$ws.Range("A134").Value = '''Q26'', '
$ws.Range("A133").Value = '''Q25'','

$ws.Range("A119:A138").Select()
